# Generate Report for Handback
#
# 1) The "Ready for handoff" status becomes "Handed back: in sync with en-US"
#    everywhere it is used (Overview!B2:C2/B3:C3, zh-cn!B2/B3, de-de!B2/B3).
# 2) Each localized sheet (zh-cn, de-de) gains two new populated columns:
#       E = Latest Target File   (hyperlink to the source .md, same as column A)
#       F = Latest Handback File (hyperlink to the handed-back .xlf, same as column C)
#    for the two real content rows (2 and 3).
# 3) The "Latest Handback DateTime" column (G) is stamped with the real
#    handback timestamp for rows 2 and 3 (previously the 0001-01-01 placeholder).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- 1) Status text, every sheet that shows it -----------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

# --- 2) & 3) Per-locale handback details ------------------------------------

$sourceRepo = "https://github.com/OpenLocalizationTest/oltest/blob/8874d8af3bed1618e54dc615f54298ba6a91a86b/e2e"

function Fill-Handback {
    param($ws, $row2Name, $row2Target, $row2TargetUrl, $row3Name, $row3Target, $row3TargetUrl, $handbackStamp)

    # Row 2 - Latest Target File / Latest Handback File
    $ws.Hyperlinks.Add($ws.Range("E2"), "$sourceRepo/$row2Name", "", "", $row2Name) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $row2TargetUrl, "", "", $row2Target) | Out-Null
    $ws.Range("G2").Value = $handbackStamp

    # Row 3 - Latest Target File / Latest Handback File
    $ws.Hyperlinks.Add($ws.Range("E3"), "$sourceRepo/$row3Name", "", "", $row3Name) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $row3TargetUrl, "", "", $row3Target) | Out-Null
    $ws.Range("G3").Value = $handbackStamp
}

Fill-Handback $wsZh `
    "35dcb9aa-cc29-4195-8928-a939010c5a18.md" `
    "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/48fdfed748c2a9e42d55b40e8068ceca9457ac3c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.zh-cn.xlf" `
    "aa278856-594c-4cac-a891-ca314d641da9.md" `
    "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/48fdfed748c2a9e42d55b40e8068ceca9457ac3c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.zh-cn.xlf" `
    "2016-02-18 08:27:59"

Fill-Handback $wsDe `
    "35dcb9aa-cc29-4195-8928-a939010c5a18.md" `
    "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4762c7b7aa50fded36d6f5f1ea9ab0718baf6ca2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.de-de.xlf" `
    "aa278856-594c-4cac-a891-ca314d641da9.md" `
    "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4762c7b7aa50fded36d6f5f1ea9ab0718baf6ca2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.de-de.xlf" `
    "2016-02-18 08:28:21"
